$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E (rows 3:7) did not previously carry the shared "left aligned"
# style (s="1") that all the other surrounding data columns use. Apply it
# so that, once the values are cleared below, the cell keeps the same
# formatting as its neighbours (and matches rows 8 and below).
$ws.Range("E3:E7").HorizontalAlignment = -4131

# Clear out the trial configuration rows (3-7): this removes the values
# (and the A-column helper formulas / D-column values entirely, since
# those columns have no special style to retain) while leaving the
# existing cell formatting (style s="1") on B, C, E:N untouched.
$ws.Range("A3:N7").ClearContents()

# Reflect the last selected cell as recorded in the saved workbook.
$ws.Range("N6").Select()
